{"js": "// Fix subject-verb agreement: \"My grandmother are already retired.\" -> \"My grandmother is already retired.\"\nconst results = context.document.body.search(\"My grandmother are already retired.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const range = results.items[0];\n  range.insertText(\"My grandmother is already retired.\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Fix subject-verb agreement in the \"grandmother\" paragraph:\n#   \"My grandmother are already retired.\" -> \"My grandmother is already retired.\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"My grandmother are already retired.\"\n$find.Replacement.Text = \"My grandmother is already retired.\"\n\n$find.Execute(\n  $find.Text,              # FindText\n  $false,                  # MatchCase\n  $false,                  # MatchWholeWord\n  $false,                  # MatchWildcards\n  $false,                  # MatchSoundsLike\n  $false,                  # MatchAllWordForms\n  $true,                   # Forward\n  1,                       # Wrap (wdFindContinue)\n  $false,                  # Format\n  $find.Replacement.Text,  # ReplaceWith\n  2                        # Replace (wdReplaceAll)\n)\n"}
